$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("N3").Value = 7.5

# Row 6
$ws.Range("G6").Value  = 1.39
$ws.Range("H6").Value  = 4.3
$ws.Range("I6").Value  = 6.2
$ws.Range("J6").Value  = 1.83
$ws.Range("K6").Value  = 2.4
$ws.Range("L6").Value  = 5.9
$ws.Range("M6").Value  = 1.02
$ws.Range("O6").Value  = 1.14
$ws.Range("R6").Value  = 2.12
$ws.Range("T6").Value  = 3.32
$ws.Range("U6").Value  = 1.85
$ws.Range("V6").Value  = 1.91
$ws.Range("W6").Value  = 6.8
$ws.Range("X6").Value  = 6.3
$ws.Range("Y6").Value  = 7
$ws.Range("Z6").Value  = 8
$ws.Range("AA6").Value = 9
$ws.Range("AB6").Value = 18.5
$ws.Range("AD6").Value = 7.6
$ws.Range("AE6").Value = 14.5
$ws.Range("AF6").Value = 55
$ws.Range("AH6").Value = 15.5
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 16.5
$ws.Range("AK6").Value = 100
$ws.Range("AN6").Value = 3.35
$ws.Range("AO6").Value = 6.2
$ws.Range("AP6").Value = 14.5
$ws.Range("AQ6").Value = 17
$ws.Range("AU6").Value = 7.8
$ws.Range("AW6").Value = 7.9
$ws.Range("AX6").Value = 35
$ws.Range("AY6").Value = 35

# Row 8
$ws.Range("J8").Value = 2.62
